$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text (matching the source inlineStr cells) by pre-formatting as Text
# before assignment - otherwise Excel auto-converts "5.72" etc. to a number.
$textCells = @("D5", "D6", "D10", "D11", "D13", "D14", "D18", "D20", "D21", "D22", "D24", "D25", "D28", "D31", "D32", "D37", "D39", "D40", "D44", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range("D2").Value = '63.871.87'
$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").Value = '3.135.42'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '591.72'
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("D6").Value = '147.09'
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.126.00'
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +18.82%  '
$ws.Range("D11").Value = '5.72'
$ws.Range("E11").Value = '  +4.25%  '
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  +6.16%  '
$ws.Range("D14").Value = '36.14'
$ws.Range("E14").Value = '  +3.43%  '
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = '3.653.48'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '63.798.57'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '7.17'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '3.132.80'
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("D20").Value = '466.45'
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").Value = '14.25'
$ws.Range("E21").Value = '  +2.17%  '
$ws.Range("D22").Value = '0.733'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  +3.07%  '
$ws.Range("D24").Value = '13.31'
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").Value = '82.33'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("E27").Value = '  +7.69%  '
$ws.Range("D28").Value = '2.71'
$ws.Range("E28").Value = '  +3.02%  '
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").Value = '6.86'
$ws.Range("E31").Value = '  +3.26%  '
$ws.Range("D32").Value = '27.07'
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("E33").Value = '  +2.03%  '
$ws.Range("D34").Value = '0.0₃0870'
$ws.Range("E34").Value = '  +7.47%  '
$ws.Range("E35").Value = '  +9.05%  '
$ws.Range("E36").Value = '  +2.27%  '
$ws.Range("D37").Value = '3.44'
$ws.Range("E37").Value = '  +16.06%  '
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Value = '456.01'
$ws.Range("E39").Value = '  +9.13%  '
$ws.Range("D40").Value = '50.87'
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("D43").Value = '2.920.92'
$ws.Range("E43").Value = '  +5.27%  '
$ws.Range("D44").Value = '0.277'
$ws.Range("E44").Value = '  +4.48%  '
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("E46").Value = '  +3.38%  '
$ws.Range("D47").Value = '129.51'
$ws.Range("E47").Value = '  +4.99%  '
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '34.53'
$ws.Range("E49").Value = '  -5.72%  '
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = '24.72'
$ws.Range("E51").Value = '  +2.55%  '
